$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sat at the end of the
#    "Un sistema de emision de boletos..." paragraph - it will be
#    re-created below, at the point of the newest edit.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Add the discussion answer as a new paragraph right after
#    "Ud. acepta un empleo..." and move the _GoBack bookmark there,
#    marking it as the most recently edited spot.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("también tiene una responsabilidad de confidencialidad con su empleador anterior.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$srcPara = $rng.Paragraphs(1)
$srcPara.Range.InsertParagraphAfter()
$newPara = $srcPara.Next()

$newRange = $newPara.Range
$newRange.ParagraphFormat.Style = "Normal"
$newRange.ParagraphFormat.LeftIndent = 1416 / 1440 * 72
$newRange.ParagraphFormat.Alignment = 3
$newRange.Text = "Corregiría las ambigüedades en la interpretación de los requerimientos sin plantearlo como correcciones sino como mejoras, ya que existe una responsabilidad asumida con el actual empleador, y al solucionar dicho inconveniente se lograría que la empresa o institución pueda alcanzar sus objetivos."

$insStart = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $insStart)

# ---------------------------------------------------------------------
# 3) The extra page pushes the footer's cached PAGE field result from
#    7 to 8.
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
foreach ($f in $ftr.Range.Fields) {
    if ($f.Code.Text.Trim() -eq "PAGE") {
        $f.Result.Find.Execute("7", $false, $false, $false, $false, $false, $true, 1, $false, "8", 2)
    }
}
